# Apply updated TPM values to Agrn-Dag1 sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("G2").Value = 8.382531
$ws.Range("H2").Value = 25.147593
$ws.Range("I2").Value = 0.3278601051951505
$ws.Range("J2").Value = 0.3278601051951506
$ws.Range("M2").Value = 11.319211
$ws.Range("N2").Value = 33.957633
$ws.Range("O2").Value = 0.09922284194232082
$ws.Range("P2").Value = 0.09922284194232082
$ws.Range("Q2").Value = 94.88363710304101
$ws.Range("R2").Value = 853.9527339273691
$ws.Range("S2").Value = 0.0325312113969711
$ws.Range("T2").Value = 0.03253121139697111

# Row 3
$ws.Range("G3").Value = 8.382531
$ws.Range("H3").Value = 25.147593
$ws.Range("I3").Value = 0.3278601051951505
$ws.Range("J3").Value = 0.3278601051951506
$ws.Range("O3").Value = 0.3843080175847637
$ws.Range("P3").Value = 0.3843080175847637
$ws.Range("Q3").Value = 367.501492221912
$ws.Range("R3").Value = 3307.513429997208
$ws.Range("S3").Value = 0.1259992670726804
$ws.Range("T3").Value = 0.1259992670726804

# Row 4
$ws.Range("G4").Value = 8.382531
$ws.Range("H4").Value = 25.147593
$ws.Range("I4").Value = 0.3278601051951505
$ws.Range("J4").Value = 0.3278601051951506
$ws.Range("O4").Value = 0.5164691404729155
$ws.Range("P4").Value = 0.5164691404729155
$ws.Range("Q4").Value = 493.8829561850121
$ws.Range("R4").Value = 4444.946605665109
$ws.Range("S4").Value = 0.169329626725499
$ws.Range("T4").Value = 0.1693296267254991

# Row 5
$ws.Range("I5").Value = 0.2503004183517279
$ws.Range("J5").Value = 0.250300418351728
$ws.Range("M5").Value = 11.319211
$ws.Range("N5").Value = 33.957633
$ws.Range("O5").Value = 0.09922284194232082
$ws.Range("P5").Value = 0.09922284194232082
$ws.Range("Q5").Value = 72.43764546311134
$ws.Range("R5").Value = 651.938809168002
$ws.Range("S5").Value = 0.02483551884821028
$ws.Range("T5").Value = 0.02483551884821029

# Row 6
$ws.Range("I6").Value = 0.2503004183517279
$ws.Range("J6").Value = 0.250300418351728
$ws.Range("O6").Value = 0.3843080175847637
$ws.Range("P6").Value = 0.3843080175847637
$ws.Range("S6").Value = 0.09619245757738958
$ws.Range("T6").Value = 0.09619245757738959

# Row 7
$ws.Range("I7").Value = 0.2503004183517279
$ws.Range("J7").Value = 0.250300418351728
$ws.Range("O7").Value = 0.5164691404729155
$ws.Range("P7").Value = 0.5164691404729155
$ws.Range("S7").Value = 0.1292724419261281
$ws.Range("T7").Value = 0.1292724419261281

# Row 8
$ws.Range("I8").Value = 0.4218394764531215
$ws.Range("J8").Value = 0.4218394764531215
$ws.Range("M8").Value = 11.319211
$ws.Range("N8").Value = 33.957633
$ws.Range("O8").Value = 0.09922284194232082
$ws.Range("P8").Value = 0.09922284194232082
$ws.Range("Q8").Value = 122.0815316205993
$ws.Range("R8").Value = 1098.733784585394
$ws.Range("S8").Value = 0.04185611169713944
$ws.Range("T8").Value = 0.04185611169713944

# Row 9
$ws.Range("I9").Value = 0.4218394764531215
$ws.Range("J9").Value = 0.4218394764531215
$ws.Range("O9").Value = 0.3843080175847637
$ws.Range("P9").Value = 0.3843080175847637
$ws.Range("S9").Value = 0.1621162929346937
$ws.Range("T9").Value = 0.1621162929346937

# Row 10
$ws.Range("I10").Value = 0.4218394764531215
$ws.Range("J10").Value = 0.4218394764531215
$ws.Range("O10").Value = 0.5164691404729155
$ws.Range("P10").Value = 0.5164691404729155
$ws.Range("S10").Value = 0.2178670718212883
$ws.Range("T10").Value = 0.2178670718212883
